$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell for 05-10-2020 (column T), styled like the preceding
#     date headers (bold, centered, thin border, text rather than a date
#     serial so it displays exactly as "05-10-2020") ---
$header = $ws.Range("T1")
$header.NumberFormat = "@"
$header.Value = "05-10-2020"
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous

# --- New daily totals for 05-10-2020, one per state/UT row (rows 2-36) ---
$values = @{
  2  = 3649
  3  = 658875
  4  = 7577
  5  = 152127
  6  = 175458
  7  = 10598
  8  = 93731
  9  = 2980
  10 = 260350
  11 = 30033
  12 = 122233
  13 = 120341
  14 = 12361
  15 = 62404
  16 = 75531
  17 = 515782
  18 = 144471
  19 = 3354
  20 = 113832
  21 = 1149603
  22 = 9205
  23 = 4393
  24 = 1807
  25 = 5309
  26 = 202302
  27 = 23763
  28 = 100977
  29 = 121331
  30 = 2480
  31 = 564092
  32 = 172388
  33 = 21876
  34 = 41740
  35 = 362052
  36 = 237698
}

foreach ($row in $values.Keys) {
  $ws.Cells.Item($row, 20).Value = $values[$row]
}
